$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 1468.6666
$ws.Range("I5").Value = 1162.4
$ws.Range("K5").Value = 1162.4
$ws.Range("M5").Value = -1047.4

$ws.Range("H6").Value = 43.714287
$ws.Range("I6").Value = 44.333332
$ws.Range("K6").Value = 132.999996
$ws.Range("M6").Value = -20.99999600000001

$ws.Range("H19").Value = 1639.64
$ws.Range("I19").Value = 908.5
$ws.Range("J19").Value = 2939.4443
$ws.Range("K19").Value = 908.5
$ws.Range("L19").Value = 2939.4443
$ws.Range("M19").Value = -733.5
$ws.Range("N19").Value = -3289.4443

$ws.Range("H113").Value = 7291.522
$ws.Range("I113").Value = 4399.9
$ws.Range("J113").Value = 9515.846
$ws.Range("K113").Value = 4399.9
$ws.Range("L113").Value = 9515.846
$ws.Range("M113").Value = -1145.9
$ws.Range("N113").Value = -16023.846

$ws.Range("H137").Value = 3029.6943
$ws.Range("I137").Value = 2504.5
$ws.Range("J137").Value = 4867.875
$ws.Range("K137").Value = 7513.5
$ws.Range("L137").Value = 14603.625
$ws.Range("M137").Value = -4963.5
$ws.Range("N137").Value = -19703.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1839.5167
$ws.Range("I32").Value = 1068
$ws.Range("K32").Value = 1068
$ws.Range("M32").Value = -781

$ws.Range("H41").Value = 7367
$ws.Range("I41").Value = 7367
$ws.Range("K41").Value = 7367
$ws.Range("M41").Value = -6953

$ws.Range("H61").Value = 2622.2334
$ws.Range("I61").Value = 1457.6666
$ws.Range("K61").Value = 1457.6666
$ws.Range("M61").Value = -1245.6666

$ws.Range("H74").Value = 2530.1482
$ws.Range("I74").Value = 1625.2
$ws.Range("J74").Value = 3310.276
$ws.Range("K74").Value = 1625.2
$ws.Range("L74").Value = 3310.276
$ws.Range("M74").Value = -751.2
$ws.Range("N74").Value = -5058.276

$ws.Range("H77").Value = 2530.1482
$ws.Range("I77").Value = 1625.2
$ws.Range("J77").Value = 3310.276
$ws.Range("K77").Value = 8126
$ws.Range("L77").Value = 16551.38
$ws.Range("M77").Value = -3758
$ws.Range("N77").Value = -25287.38

$ws.Range("H111").Value = 99749.5
$ws.Range("J111").Value = 99749.5
$ws.Range("L111").Value = 99749.5
$ws.Range("N111").Value = -107929.5

$ws.Range("H119").Value = 49600
$ws.Range("J119").Value = 49600
$ws.Range("L119").Value = 49600
$ws.Range("N119").Value = -59276

$ws.Range("H132").Value = 2593.5247
$ws.Range("I132").Value = 2417.6
$ws.Range("K132").Value = 7252.799999999999
$ws.Range("M132").Value = -4722.799999999999

$ws.Range("H136").Value = 2622.2334
$ws.Range("I136").Value = 1457.6666
$ws.Range("K136").Value = 4372.9998
$ws.Range("M136").Value = -1822.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H14").Value = 7999
$ws.Range("J14").Value = 4009
$ws.Range("L14").Value = 4009
$ws.Range("N14").Value = -4353

$ws.Range("H107").Value = 1863.091
$ws.Range("I107").Value = 1724.4
$ws.Range("K107").Value = 1724.4
$ws.Range("M107").Value = 195.5999999999999

$ws.Range("H134").Value = 2271.4333
$ws.Range("I134").Value = 2145.5186
$ws.Range("J134").Value = 3404.6667
$ws.Range("K134").Value = 6436.5558
$ws.Range("L134").Value = 10214.0001
$ws.Range("M134").Value = -3901.5558
$ws.Range("N134").Value = -15284.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 539.5
$ws.Range("I14").Value = 380
$ws.Range("J14").Value = 699
$ws.Range("K14").Value = 380
$ws.Range("L14").Value = 699
$ws.Range("M14").Value = -210
$ws.Range("N14").Value = -1039

$ws.Range("H16").Value = 1916
$ws.Range("I16").Value = 1970.3684
$ws.Range("K16").Value = 1970.3684
$ws.Range("M16").Value = -1683.3684

$ws.Range("H88").Value = 13297
$ws.Range("J88").Value = 13297
$ws.Range("L88").Value = 13297
$ws.Range("N88").Value = -14109

$ws.Range("H91").Value = 13297
$ws.Range("J91").Value = 13297
$ws.Range("L91").Value = 13297
$ws.Range("N91").Value = -16105

$ws.Range("H113").Value = 1916
$ws.Range("I113").Value = 1970.3684
$ws.Range("K113").Value = 1970.3684
$ws.Range("M113").Value = 199.6315999999999

$ws.Range("H134").Value = 3548.3845
$ws.Range("I134").Value = 3437.5938
$ws.Range("K134").Value = 10312.7814
$ws.Range("M134").Value = -7777.7814

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 118.57143
$ws.Range("J38").Value = 121.14286
$ws.Range("L38").Value = 363.42858
$ws.Range("N38").Value = -1057.42858

$ws.Range("H68").Value = 1001.5
$ws.Range("I68").Value = 1002
$ws.Range("J68").Value = 1001
$ws.Range("K68").Value = 3006
$ws.Range("L68").Value = 3003
$ws.Range("M68").Value = -2195
$ws.Range("N68").Value = -4625

$ws.Range("H71").Value = 1001.5
$ws.Range("I71").Value = 1002
$ws.Range("J71").Value = 1001
$ws.Range("K71").Value = 9018
$ws.Range("L71").Value = 9009
$ws.Range("M71").Value = -4962
$ws.Range("N71").Value = -17121

$ws.Range("H116").Value = 1621.3334
$ws.Range("I116").Value = 1405.6
$ws.Range("K116").Value = 4216.799999999999
$ws.Range("M116").Value = -774.7999999999993

$ws.Range("H134").Value = 1518.1666
$ws.Range("I134").Value = 1518.1666
$ws.Range("K134").Value = 4554.4998
$ws.Range("M134").Value = 515.5002000000004

$ws.Range("H140").Value = 5545.846
$ws.Range("I140").Value = 1209.6
$ws.Range("K140").Value = 3628.8
$ws.Range("M140").Value = 1551.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 42267.43
$ws.Range("I70").Value = 62141.11
$ws.Range("K70").Value = 62141.11
$ws.Range("M70").Value = -61871.11

$ws.Range("H73").Value = 42267.43
$ws.Range("I73").Value = 62141.11
$ws.Range("K73").Value = 62141.11
$ws.Range("M73").Value = -61205.11

$ws.Range("H92").Value = 29998.5
$ws.Range("J92").Value = 29998.5
$ws.Range("L92").Value = 29998.5
$ws.Range("N92").Value = -33742.5

$ws.Range("H134").Value = 12995
$ws.Range("J134").Value = 12995
$ws.Range("L134").Value = 38985
$ws.Range("N134").Value = -44055

$ws.Range("H136").Value = 29550.334
$ws.Range("J136").Value = 29550.334
$ws.Range("L136").Value = 88651.00199999999
$ws.Range("N136").Value = -93751.00199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3881.8
$ws.Range("I40").Value = 4013.647
$ws.Range("K40").Value = 4013.647
$ws.Range("M40").Value = -3877.647

$ws.Range("H61").Value = 8592.362999999999
$ws.Range("I61").Value = 8961.6
$ws.Range("K61").Value = 8961.6
$ws.Range("M61").Value = -8759.6

$ws.Range("H113").Value = 8592.362999999999
$ws.Range("I113").Value = 8961.6
$ws.Range("K113").Value = 8961.6
$ws.Range("M113").Value = -6791.6

$ws.Range("H136").Value = 3336.9412
$ws.Range("I136").Value = 2871
$ws.Range("K136").Value = 8613
$ws.Range("M136").Value = -6063

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3395688.8
$ws.Range("I96").Value = 3734257.8
$ws.Range("K96").Value = 3734257.8
$ws.Range("M96").Value = -3732884.8

$ws.Range("H122").Value = 4147.593
$ws.Range("I122").Value = 4336.75
$ws.Range("J122").Value = 3607.1428
$ws.Range("K122").Value = 13010.25
$ws.Range("L122").Value = 10821.4284
$ws.Range("M122").Value = -10560.25
$ws.Range("N122").Value = -15721.4284

$ws.Range("H136").Value = 6715.9165
$ws.Range("I136").Value = 7629.1
$ws.Range("K136").Value = 22887.3
$ws.Range("M136").Value = -20337.3
